# Updates coin prices / 1h volume percentages to the latest scraped values,
# and swaps the Fetch.AI / Stacks rows (ranking order changed upstream).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells that contain plain decimal numbers (e.g. "9.00", "0.610")
# must be forced to Text format first, otherwise Excel would silently
# convert them to numeric values and drop significant trailing zeros.
$numericPriceCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D16", "D20", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D40", "D45", "D47", "D51", "D41", "D42")
foreach ($cellRef in $numericPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.412.82'
$ws.Range("E2").Value = '  -2.78%  '
$ws.Range("D3").Value = '3.436.49'
$ws.Range("E3").Value = '  -5.09%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '571.38'
$ws.Range("E5").Value = '  -5.06%  '
$ws.Range("D6").Value = '188.98'
$ws.Range("E6").Value = '  -3.85%  '
$ws.Range("D7").Value = '0.603'
$ws.Range("E7").Value = '  -3.78%  '
$ws.Range("D8").Value = '3.425.71'
$ws.Range("E8").Value = '  -5.09%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").Value = '0.201'
$ws.Range("E10").Value = '  -5.42%  '
$ws.Range("D11").Value = '0.610'
$ws.Range("E11").Value = '  -5.55%  '
$ws.Range("D12").Value = '50.75'
$ws.Range("E12").Value = '  -4.65%  '
$ws.Range("D13").Value = '0.0000281'
$ws.Range("E13").Value = '  -7.47%  '
$ws.Range("D14").Value = '9.00'
$ws.Range("E14").Value = '  -5.72%  '
$ws.Range("D15").Value = '3.987.85'
$ws.Range("E15").Value = '  -4.96%  '
$ws.Range("D16").Value = '630.70'
$ws.Range("E16").Value = '  +4.08%  '
$ws.Range("D17").Value = '68.395.83'
$ws.Range("E17").Value = '  -2.92%  '
$ws.Range("D18").Value = '3.442.57'
$ws.Range("E18").Value = '  -5.30%  '
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("D20").Value = '12.19'
$ws.Range("E20").Value = '  -5.54%  '
$ws.Range("E21").Value = '  -5.53%  '
$ws.Range("D22").Value = '0.932'
$ws.Range("E22").Value = '  -6.61%  '
$ws.Range("D23").Value = '17.68'
$ws.Range("E23").Value = '  -2.62%  '
$ws.Range("D24").Value = '5.31'
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("D25").Value = '98.49'
$ws.Range("E25").Value = '  -4.55%  '
$ws.Range("D26").Value = '4.22'
$ws.Range("E26").Value = '  -8.28%  '
$ws.Range("E27").Value = '  -5.96%  '
$ws.Range("D28").Value = '6.07'
$ws.Range("E28").Value = '  +1.97%  '
$ws.Range("D29").Value = '9.76'
$ws.Range("E29").Value = '  -8.16%  '
$ws.Range("D30").Value = '9.08'
$ws.Range("E30").Value = '  -6.19%  '
$ws.Range("D31").Value = '32.06'
$ws.Range("E31").Value = '  -5.14%  '
$ws.Range("D32").Value = '4.10'
$ws.Range("E32").Value = '  -12.28%  '
$ws.Range("D33").Value = '6.63'
$ws.Range("E33").Value = '  -8.82%  '
$ws.Range("D34").Value = '11.47'
$ws.Range("E34").Value = '  -6.50%  '
$ws.Range("D35").Value = '60.49'
$ws.Range("E35").Value = '  -4.48%  '
$ws.Range("E36").Value = '  -7.68%  '
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").Value = '3.623.84'
$ws.Range("E38").Value = '  -7.49%  '
$ws.Range("D39").Value = '0.0₃0776'
$ws.Range("E39").Value = '  -12.72%  '
$ws.Range("D40").Value = '500.48'
$ws.Range("E40").Value = '  -4.01%  '
$ws.Range("E43").Value = '  -6.25%  '
$ws.Range("E44").Value = '  -3.20%  '
$ws.Range("D45").Value = '33.95'
$ws.Range("E45").Value = '  -7.53%  '
$ws.Range("E46").Value = '  +64.80%  '
$ws.Range("D47").Value = '0.0434'
$ws.Range("E47").Value = '  -5.88%  '
$ws.Range("E48").Value = '  -4.31%  '
$ws.Range("E49").Value = '  -4.09%  '
$ws.Range("E50").Value = '  -4.95%  '
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.42%  '

# Row 41/42: Fetch.AI moves up to rank 41 (was Stacks), Stacks drops to 42,
# each carrying its own refreshed price / volume figures.
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = '2.87'
$ws.Range("E41").Value = '  -6.31%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '3.44'
$ws.Range("E42").Value = '  -2.63%  '

Write-Host "Applied crypto list update"
